$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to the generic "Sheet1"
$ws.Name = "Sheet1"

# Fix the Lens "Angle between side rays and surface (degrees)" value
# (was "71.33/80.08", corrected to "67.15/88.13")
$ws.Range("C8").Value = "67.15/88.13"

# Remove the old tolerance/finesse block (rows 20-24: Right NA, Left NA,
# Finesse, roundtrip power losses, Free Spectral Range) - these stats are
# no longer kept on this sheet
$ws.Rows("20:24").Delete()

# Update the selected cell shown when the workbook is reopened
$ws.Range("N3").Select() | Out-Null
